# "update scripts wuth new tpm" - recomputed NATMI ligand/receptor-derived
# metrics (Jag2/Notch2) after the sending-cluster ligand average/total and
# target-cluster receptor average/total expression values were refreshed
# with new TPM figures (cluster "ECs" for both Jag2-sending and
# Notch2-receiving). Every dependent column (detection specificity, edge
# weights, edge specificity) is updated in lock-step with the new inputs.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 19.35876433333334
$ws.Range("H2").Value = 58.07629300000001
$ws.Range("I2").Value = 0.9707667559429034
$ws.Range("J2").Value = 0.9707667559429034
$ws.Range("M2").Value = 1.492477333333333
$ws.Range("N2").Value = 4.477432
$ws.Range("O2").Value = 0.02769484181536182
$ws.Range("P2").Value = 0.02769484181536182
$ws.Range("Q2").Value = 28.89251696884179
$ws.Range("R2").Value = 260.032652719576
$ws.Range("S2").Value = 0.02688523174545067
$ws.Range("T2").Value = 0.02688523174545067
$ws.Range("G3").Value = 19.35876433333334
$ws.Range("H3").Value = 58.07629300000001
$ws.Range("I3").Value = 0.9707667559429034
$ws.Range("J3").Value = 0.9707667559429034
$ws.Range("O3").Value = 0.6282762845978157
$ws.Range("P3").Value = 0.6282762845978156
$ws.Range("Q3").Value = 655.4463583826794
$ws.Range("R3").Value = 5899.017225444113
$ws.Range("S3").Value = 0.6099097306348819
$ws.Range("T3").Value = 0.6099097306348817
$ws.Range("G4").Value = 19.35876433333334
$ws.Range("H4").Value = 58.07629300000001
$ws.Range("I4").Value = 0.9707667559429034
$ws.Range("J4").Value = 0.9707667559429034
$ws.Range("N4").Value = 55.619234
$ws.Range("O4").Value = 0.3440288735868225
$ws.Range("P4").Value = 0.3440288735868225
$ws.Range("Q4").Value = 358.9065478021736
$ws.Range("R4").Value = 3230.158930219562
$ws.Range("S4").Value = 0.333971793562571
$ws.Range("T4").Value = 0.333971793562571
$ws.Range("H5").Value = 0.9049070000000001
$ws.Range("I5").Value = 0.01512585579145048
$ws.Range("J5").Value = 0.01512585579145048
$ws.Range("M5").Value = 1.492477333333333
$ws.Range("N5").Value = 4.477432
$ws.Range("O5").Value = 0.02769484181536182
$ws.Range("P5").Value = 0.02769484181536182
$ws.Range("Q5").Value = 0.450184395424889
$ws.Range("R5").Value = 4.051659558824001
$ws.Range("S5").Value = 0.0004189081834661955
$ws.Range("T5").Value = 0.0004189081834661955
$ws.Range("H6").Value = 0.9049070000000001
$ws.Range("I6").Value = 0.01512585579145048
$ws.Range("J6").Value = 0.01512585579145048
$ws.Range("O6").Value = 0.6282762845978157
$ws.Range("P6").Value = 0.6282762845978156
$ws.Range("R6").Value = 91.91464717668802
$ws.Range("S6").Value = 0.009503216478014859
$ws.Range("T6").Value = 0.009503216478014859
$ws.Range("H7").Value = 0.9049070000000001
$ws.Range("I7").Value = 0.01512585579145048
$ws.Range("J7").Value = 0.01512585579145048
$ws.Range("N7").Value = 55.619234
$ws.Range("O7").Value = 0.3440288735868225
$ws.Range("P7").Value = 0.3440288735868225
$ws.Range("S7").Value = 0.005203731129969424
$ws.Range("T7").Value = 0.005203731129969425
$ws.Range("G8").Value = 0.2813256666666666
$ws.Range("H8").Value = 0.843977
$ws.Range("I8").Value = 0.01410738826564608
$ws.Range("J8").Value = 0.01410738826564608
$ws.Range("M8").Value = 1.492477333333333
$ws.Range("N8").Value = 4.477432
$ws.Range("O8").Value = 0.02769484181536182
$ws.Range("P8").Value = 0.02769484181536182
$ws.Range("Q8").Value = 0.4198721807848889
$ws.Range("R8").Value = 3.778849627064
$ws.Range("S8").Value = 0.0003907018864449598
$ws.Range("T8").Value = 0.0003907018864449598
$ws.Range("G9").Value = 0.2813256666666666
$ws.Range("H9").Value = 0.843977
$ws.Range("I9").Value = 0.01410738826564608
$ws.Range("J9").Value = 0.01410738826564608
$ws.Range("O9").Value = 0.6282762845978157
$ws.Range("P9").Value = 0.6282762845978156
$ws.Range("Q9").Value = 9.525085411507556
$ws.Range("R9").Value = 85.72576870356801
$ws.Range("S9").Value = 0.008863337484918942
$ws.Range("T9").Value = 0.008863337484918942
$ws.Range("G10").Value = 0.2813256666666666
$ws.Range("H10").Value = 0.843977
$ws.Range("I10").Value = 0.01410738826564608
$ws.Range("J10").Value = 0.01410738826564608
$ws.Range("N10").Value = 55.619234
$ws.Range("O10").Value = 0.3440288735868225
$ws.Range("P10").Value = 0.3440288735868225
$ws.Range("R10").Value = 46.94135425361799
$ws.Range("S10").Value = 0.004853348894282179
$ws.Range("T10").Value = 0.00485334889428218
